$d = $word.ActiveDocument

# The log table is the first (only) table in the document.
$t = $d.Tables.Item(1)

# Add a new row at the end of the table and fill in the two cells,
# mirroring the existing rows (date in col 1, bulleted objective in col 2).
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "13/08/2020"
$newRow.Cells.Item(2).Range.Text = "Mejora ER"
